# "updated INTC AMD models"
#
# HardwareSemis.xlsx - Semiconductors sheet: refresh the cached AMD (row 11)
# and Intel (row 12) model outputs, fix a mislabeled quarter tag, and move
# the saved cursor position. These pull their raw inputs from external
# "AMD.xlsx" / "INTC.xlsx" workbooks ([3]/[4]) that aren't available in this
# session, so the new AMD/Intel input numbers are entered directly (the
# dependent formulas -- E11/G11/E12/G12 -- stay live formulas and recompute
# off of them automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semiconductors")

# --- AMD (row 11) ---
# Updated share count input; E11 (=D11*H11) and G11 (=E11-F11) recalc on their own.
$ws.Range("D11").Value = 130
# "as of" date pushed out to the new model refresh date.
$ws.Range("J11").Value = 45507

# --- Intel / INTC (row 12) ---
# Updated share count input, now formatted to match the AMD row (2dp, right aligned).
$ws.Range("D12").NumberFormat = "#,##0.00"
$ws.Range("D12").Value = 21
# Refreshed external-model figures (normally pulled live from [4]Main!$L$3 and
# $L$5-$L$6 in INTC.xlsx); entered directly since that workbook isn't reachable
# here. E12/G12 formulas are untouched and recompute from these.
$ws.Range("H12").Value = 4267
$ws.Range("F12").Value = -17932
# Quarter label typo fix: "Q422" -> "Q224".
$ws.Range("I12").Value = "Q224"
# "as of" date pushed out to the new model refresh date.
$ws.Range("J12").Value = 45507

# Restore the saved cursor/selection on the frozen pane.
$ws.Range("L17").Select() | Out-Null
